$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (values 2 through 36) from 1.64 to 40
$ws.Range("C2:C36").Value = 40

# Update the selection to C2:C36 with active cell C2
$ws.Range("C2:C36").Select()
